$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the LV/EXP/maxEXP/startItemID/startInventory headers to their
# underscore-prefixed equivalents (S1:W1) as part of the "exp 동기화" work.
$ws.Range("S1").Value = "_lv"
$ws.Range("T1").Value = "_exp"
$ws.Range("U1").Value = "_maxEXP"
$ws.Range("V1").Value = "_startItemID"
$ws.Range("W1").Value = "_startInventory"

# Move the active selection to Y11, matching the cursor position left in
# the sheet after the options-UI edits.
$ws.Range("Y11").Select()
